$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.420.67"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "3.429.90"
$ws.Range("E3").Value = "  -1.15%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'407.22"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").Value = "'134.79"
$ws.Range("E6").Value = "  +4.23%  "

$ws.Range("E7").Value = "  -0.67%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.687"
$ws.Range("E9").Value = "  -0.67%  "

$ws.Range("D10").Value = "'0.122"
$ws.Range("E10").Value = "  -3.26%  "

$ws.Range("D11").Value = "'42.26"
$ws.Range("E11").Value = "  -0.82%  "

$ws.Range("E12").Value = "  -0.80%  "

$ws.Range("D13").Value = "'8.47"
$ws.Range("E13").Value = "  -2.51%  "

$ws.Range("D14").Value = "'19.92"
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").Value = "3.413.30"
$ws.Range("E15").Value = "  -1.72%  "

$ws.Range("D16").Value = "62.399.39"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "'11.36"
$ws.Range("E17").Value = "  +3.86%  "

$ws.Range("E18").Value = "  -2.05%  "

$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("D20").Value = "'3.19"
$ws.Range("E20").Value = "  -4.69%  "

$ws.Range("D21").Value = "'84.32"
$ws.Range("E21").Value = "  +2.29%  "

$ws.Range("D22").Value = "'314.43"
$ws.Range("E22").Value = "  +1.56%  "

$ws.Range("D23").Value = "'12.96"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").Value = "'3.15"
$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("D25").Value = "'4.76"
$ws.Range("E25").Value = "  +9.08%  "

$ws.Range("D26").Value = "'29.77"
$ws.Range("E26").Value = "  -1.81%  "

$ws.Range("D27").Value = "'8.22"
$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("E28").Value = "  +3.49%  "

$ws.Range("D29").Value = "'7.56"
$ws.Range("E29").Value = "  -2.91%  "

$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("E31").Value = "  -2.90%  "

$ws.Range("D32").Value = "'42.33"

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "'11.37"
$ws.Range("E34").Value = "  -4.58%  "

$ws.Range("E35").Value = "  -1.37%  "

$ws.Range("D36").Value = "'51.44"
$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "'3.40"
$ws.Range("E38").Value = "  -4.69%  "

$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("E40").Value = "  +9.26%  "

$ws.Range("D41").Value = "'137.88"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.98"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.125"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").Value = "'4.04"
$ws.Range("E44").Value = "  +2.13%  "

$ws.Range("D45").Value = "'16.80"
$ws.Range("E45").Value = "  -3.97%  "

$ws.Range("E46").Value = "  -1.16%  "

$ws.Range("D47").Value = "'21.39"
$ws.Range("E47").Value = "  -4.22%  "

$ws.Range("D48").Value = "2.121.16"
$ws.Range("E48").Value = "  -3.99%  "

$ws.Range("E49").Value = "  -4.08%  "

$ws.Range("D50").Value = "'1.91"
$ws.Range("E50").Value = "  +2.92%  "

$ws.Range("D51").Value = "'1.69"
$ws.Range("E51").Value = "  +20.28%  "
